# Weekly update: insert a new price observation as row 468, pushing the
# existing rows (468:552) down by one and extending the table to row 553.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 468 (shifts 468:552 -> 469:553,
# carries the formatting - incl. the date style on column D - down from row 467).
$ws.Rows("468:468").Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Range("A468").Value = 4
$ws.Range("B468").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C468").Value = "Los Lagos"
$ws.Range("D468").Value = 45258
$ws.Range("E468").Value = 10
$ws.Range("F468").Value = 100112040
$ws.Range("G468").Value = "Cilantro"
$ws.Range("H468").Value = "Sin especificar"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 100
$ws.Range("K468").Value = 13000
$ws.Range("L468").Value = 15000
$ws.Range("M468").Value = 14000
$ws.Range("N468").Value = "$/docena de atados (2 kilos)"
$ws.Range("O468").Value = "Región de La Araucanía"
$ws.Range("P468").Value = 7000
$ws.Range("Q468").Value = 2
$ws.Range("R468").Value = "Hortaliza"
